$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B136").Value = 48654
$ws.Range("E136").Value = 38.26
$ws.Range("F136").Value = -1
$ws.Range("G136").Value = -32.02
$ws.Range("B137").Value = 63902
$ws.Range("E137").Value = 34.04
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("B146").Value = 53925
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 66.44
$ws.Range("B147").Value = 64350
$ws.Range("E147").Value = 70.63
$ws.Range("F147").Value = 2
$ws.Range("G147").Value = 132.88
$ws.Range("B148").Value = 57756
$ws.Range("E148").Value = 79.37
$ws.Range("F148").Value = -100
$ws.Range("G148").Value = -6644
$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 5
$ws.Range("G246").Value = 166.5
$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2
$ws.Range("B292").Value = 55373
$ws.Range("E292").Value = 163.62
$ws.Range("F292").Value = -94
$ws.Range("G292").Value = -13562.32
$ws.Range("B293").Value = 63520
$ws.Range("E293").Value = 153.4
$ws.Range("F293").Value = 72
$ws.Range("G293").Value = 10388.16
$ws.Range("B294").Value = 63571
$ws.Range("F294").Value = 0
$ws.Range("G294").Value = 0
$ws.Range("B296").Value = 63531
$ws.Range("F296").Value = 80
$ws.Range("G296").Value = 11478.4
$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12
$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 133
$ws.Range("G300").Value = 6336.12
$ws.Range("B311").Value = 61605
$ws.Range("E311").Value = 133.78
$ws.Range("F311").Value = -13
$ws.Range("G311").Value = -1455.48
$ws.Range("B312").Value = 63563
$ws.Range("E312").Value = 119.04
$ws.Range("F312").Value = 0
$ws.Range("G312").Value = 0
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 42
$ws.Range("G420").Value = 4432.68
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2
$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28
$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0
$ws.Range("B476").Value = 45706
$ws.Range("E476").Value = 23.58
$ws.Range("F476").Value = -202
$ws.Range("G476").Value = -3985.46
$ws.Range("B477").Value = 64922
$ws.Range("E477").Value = 20.98
$ws.Range("F477").Value = 68
$ws.Range("G477").Value = 1341.64
$ws.Range("B479").Value = 64927
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 131
$ws.Range("G479").Value = 2124.82
$ws.Range("B480").Value = 45718
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68
$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 145
$ws.Range("G485").Value = 1906.75
$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945
$ws.Range("B590").Value = 64833
$ws.Range("E590").Value = 34.9
$ws.Range("F590").Value = 95
$ws.Range("G590").Value = 3118.85
$ws.Range("B591").Value = 60025
$ws.Range("E591").Value = 37.22
$ws.Range("F591").Value = -98
$ws.Range("G591").Value = -3217.34
$ws.Range("B710").Value = 63150
$ws.Range("D710").Value = 75.68000000000001
$ws.Range("E710").Value = 80.45
$ws.Range("F710").Value = 32
$ws.Range("G710").Value = 2421.76
$ws.Range("B711").Value = 61428
$ws.Range("D711").Value = 69.16
$ws.Range("E711").Value = 73.52
$ws.Range("F711").Value = 1
$ws.Range("G711").Value = 69.16
$ws.Range("B737").Value = 65079
$ws.Range("F737").Value = 21
$ws.Range("G737").Value = 858.27
$ws.Range("B738").Value = 65362
$ws.Range("F738").Value = 40
$ws.Range("G738").Value = 1634.8
